# stock_tracker.xlsx edit:
#  - fill in the "clonal stock count" column (D3:D19) with the new counts
#    collected for this batch (E column already holds the =D/0.3 dilution
#    formula, so it recalculates automatically once D is populated)
#  - move the selection cursor to where the user left off (C25)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# D3:D19, in row order
$counts = @(13, 20, 18, 8, 22, 8, 25, 57, 8, 24, 23, 14, 1, 3, 11, 25, 20)

$startRow = 3
for ($i = 0; $i -lt $counts.Length; $i++) {
    $ws.Cells.Item($startRow + $i, 4).Value = $counts[$i]
}

$ws.Range("C25").Select()
